$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.244275331497192
$ws.Range("B1").Value = 2.269115447998047
$ws.Range("C1").Value = 2.954982042312622
$ws.Range("D1").Value = 3.422497749328613
$ws.Range("E1").Value = 1.742268800735474
